$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.829.97'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.58%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.357.39'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.93%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.16%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.07'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.90%  '

$ws.Range('E7').Value = '  +0.20%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.93'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.85%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.121'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.91%  '

$ws.Range('E11').Value = '  +0.81%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.932.78'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.92%  '

$ws.Range('E13').Value = '  +0.98%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.93'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.25%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.360.24'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.96%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000168'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.12%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.977.96'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.42%  '

$ws.Range('E18').Value = '  -1.61%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.15'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.20%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.83'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.87%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '374.11'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.74%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.38'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.94%  '

$ws.Range('E24').Value = '  +0.04%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.502.99'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.55%  '

$ws.Range('E26').Value = '  -6.06%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.175'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.50%  '

$ws.Range('E28').Value = '  +0.01%  '

$ws.Range('E29').Value = '  -4.08%  '

$ws.Range('E30').Value = '  +0.00%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.67'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.79%  '

$ws.Range('E33').Value = '  -2.29%  '

$ws.Range('E34').Value = '  -3.04%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.36'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.19%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '168.59'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.14%  '

$ws.Range('E37').Value = '  -5.95%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.77'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.96%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '29.22'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.58%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.392.86'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.92%  '

$ws.Range('E41').Value = '  -3.36%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.28'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.36%  '

$ws.Range('E43').Value = '  -4.30%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.28'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.81%  '

$ws.Range('E45').Value = '  -3.91%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.490.68'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.60%  '

$ws.Range('E48').Value = '  -3.21%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.53'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.38%  '

$ws.Range('E50').Value = '  -0.01%  '

$ws.Range('E51').Value = '  -2.65%  '

